$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> Username (column B), Password (column C)
$data = @{
    2  = @{ B = "jsmith";  C = "fluffy" }
    3  = @{ B = "msmith";  C = "abby" }
    4  = @{ B = "mr11j";   C = "password" }
    6  = @{ B = "test";    C = "test123" }
    7  = @{ B = "admin";   C = "admin123" }
    8  = @{ B = "llbean";  C = "blergh" }
    9  = @{ B = "firstprez"; C = "delaware" }
    10 = @{ B = "bmicro";  C = "melinda" }
    11 = @{ B = "sapple";  C = "blackturtleneck" }
    12 = @{ B = "fsusarah"; C = "spanx" }
    13 = @{ B = "model";   C = "fashion" }
    14 = @{ B = "leodicap"; C = "nooscar" }
    15 = @{ B = "heisenberg"; C = "bluesky" }
    16 = @{ B = "desanta"; C = "gta5" }
    17 = @{ B = "drake";   C = "nothingwasthesame" }
    18 = @{ B = "donglover"; C = "thedeepweb" }
    19 = @{ B = "iggy";    C = "work" }
    20 = @{ B = "god";     C = "omnicient" }
}

foreach ($row in ($data.Keys | Sort-Object)) {
    $ws.Range("B$row").Value = $data[$row].B
    $ws.Range("C$row").Value = $data[$row].C
}

# Column C width changed from 16.5 to 18 (use an input that, after the
# runtime's internal pixel-rounding, serializes to exactly width="18")
$ws.Columns.Item(3).ColumnWidth = 17.166666666666668

# Selected cell changed from B2 to C20
$ws.Range("C20").Select()
